# paises.xlsx update - "Update countries & provincias Spain"
#
# The sheet "Pais" is a COVID table (row 3 = headers, rows 4.. = one row per
# country) kept sorted descending by column B ("Casos totales"). This update
# refreshes several countries' figures; for the ones whose new total moved
# them ahead of a neighbour in the sort order, the row carrying the higher
# (new) total is written with the refreshed numbers while the row(s) below
# it keep the displaced country's previous (unchanged) figures - i.e. the
# three/two-row groups below are simple re-labelled shifts, not new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Abril de 2020 a las 14:22"

# Row 14 - Paises Bajos (in place, refreshed figures)
$ws.Cells.Item(14, 2).Value = 28153
$ws.Cells.Item(14, 3).Value = 734
$ws.Cells.Item(14, 5).Value = 24769
$ws.Cells.Item(14, 7).Value = 189
$ws.Cells.Item(14, 8).Value = 3134

# Row 20 - Austria (in place, refreshed figures)
$ws.Cells.Item(20, 2).Value = 14297
$ws.Cells.Item(20, 3).Value = 71
$ws.Cells.Item(20, 5).Value = 5806

# Rows 21-25: Israel (21, unchanged), then Suecia overtakes India/Irlanda
# Row 22 becomes Suecia with brand-new figures
$ws.Cells.Item(22, 1).Value = "Suecia"
$ws.Cells.Item(22, 2).Value = 11927
$ws.Cells.Item(22, 3).Value = 482
$ws.Cells.Item(22, 4).Value = 381
$ws.Cells.Item(22, 5).Value = 10343
$ws.Cells.Item(22, 6).Value = 954
$ws.Cells.Item(22, 7).Value = 170
$ws.Cells.Item(22, 8).Value = 1203

# Row 23 becomes India, carrying India's previous (unchanged) figures
$ws.Cells.Item(23, 1).Value = "India"
$ws.Cells.Item(23, 2).Value = 11555
$ws.Cells.Item(23, 3).Value = 68
$ws.Cells.Item(23, 4).Value = 1362
$ws.Cells.Item(23, 5).Value = 9797
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 3
$ws.Cells.Item(23, 8).Value = 396

# Row 24 becomes Irlanda, carrying Irlanda's previous (unchanged) figures
$ws.Cells.Item(24, 1).Value = "Irlanda"
$ws.Cells.Item(24, 2).Value = 11479
$ws.Cells.Item(24, 4).Value = 77
$ws.Cells.Item(24, 5).Value = 10996
$ws.Cells.Item(24, 6).Value = 194
$ws.Cells.Item(24, 8).Value = 406

# Row 33 - Dinamarca (in place, refreshed figures)
$ws.Cells.Item(33, 5).Value = 3857
$ws.Cells.Item(33, 7).Value = 10
$ws.Cells.Item(33, 8).Value = 309

# Rows 61-62: Croacia overtakes Islandia
# Row 61 becomes Croacia with brand-new figures
$ws.Cells.Item(61, 1).Value = "Croacia"
$ws.Cells.Item(61, 2).Value = 1741
$ws.Cells.Item(61, 3).Value = 37
$ws.Cells.Item(61, 4).Value = 473
$ws.Cells.Item(61, 5).Value = 1234
$ws.Cells.Item(61, 6).Value = 31
$ws.Cells.Item(61, 7).Value = 3
$ws.Cells.Item(61, 8).Value = 34

# Row 62 becomes Islandia, carrying Islandia's previous (unchanged) figures
$ws.Cells.Item(62, 1).Value = "Islandia"
$ws.Cells.Item(62, 2).Value = 1720
$ws.Cells.Item(62, 4).Value = 989
$ws.Cells.Item(62, 5).Value = 723
$ws.Cells.Item(62, 6).Value = 8
$ws.Cells.Item(62, 8).Value = 8

# Row 88 - Letonia (in place, refreshed figures)
$ws.Cells.Item(88, 4).Value = 44
$ws.Cells.Item(88, 5).Value = 617

# Rows 110-112: Senegal overtakes Estado de Palestina/Georgia
# Row 110 becomes Senegal with brand-new figures
$ws.Cells.Item(110, 1).Value = "Senegal"
$ws.Cells.Item(110, 2).Value = 314
$ws.Cells.Item(110, 3).Value = 15
$ws.Cells.Item(110, 4).Value = 190
$ws.Cells.Item(110, 5).Value = 122
$ws.Cells.Item(110, 6).Value = 1

# Row 111 becomes Estado de Palestina, carrying its previous (unchanged) figures
$ws.Cells.Item(111, 1).Value = "Estado de Palestina"
$ws.Cells.Item(111, 2).Value = 308
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = 62
$ws.Cells.Item(111, 5).Value = 244
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 8).Value = 2

# Row 112 becomes Georgia, carrying its previous (unchanged) figures
$ws.Cells.Item(112, 1).Value = "Georgia"
$ws.Cells.Item(112, 2).Value = 306
$ws.Cells.Item(112, 3).Value = 6
$ws.Cells.Item(112, 4).Value = 69
$ws.Cells.Item(112, 5).Value = 234
$ws.Cells.Item(112, 6).Value = 6
$ws.Cells.Item(112, 8).Value = 3

# Row 117 - Sri Lanka (in place, refreshed figures)
$ws.Cells.Item(117, 2).Value = 235
$ws.Cells.Item(117, 3).Value = 2
$ws.Cells.Item(117, 5).Value = 165

# Rows 118-119: Kenia overtakes Mayotte
# Row 118 becomes Kenia with brand-new figures
$ws.Cells.Item(118, 1).Value = "Kenia"
$ws.Cells.Item(118, 2).Value = 225
$ws.Cells.Item(118, 3).Value = 9
$ws.Cells.Item(118, 4).Value = 53
$ws.Cells.Item(118, 5).Value = 162
$ws.Cells.Item(118, 6).Value = 2
$ws.Cells.Item(118, 7).Value = 1
$ws.Cells.Item(118, 8).Value = 10

# Row 119 becomes Mayotte, carrying its previous (unchanged) figures
$ws.Cells.Item(119, 1).Value = "Mayotte"
$ws.Cells.Item(119, 2).Value = 217
$ws.Cells.Item(119, 4).Value = 69
$ws.Cells.Item(119, 5).Value = 145
$ws.Cells.Item(119, 6).Value = 3
$ws.Cells.Item(119, 8).Value = 3

# Rows 139-141: Togo overtakes Gabon/Liechtenstein
# Row 139 becomes Togo with brand-new figures
$ws.Cells.Item(139, 1).Value = "Togo"
$ws.Cells.Item(139, 2).Value = 81
$ws.Cells.Item(139, 3).Value = 4
$ws.Cells.Item(139, 4).Value = 35
$ws.Cells.Item(139, 5).Value = 43
$ws.Cells.Item(139, 8).Value = 3

# Row 140 becomes Gabon, carrying its previous (unchanged) figures
$ws.Cells.Item(140, 1).Value = "Gabon"
$ws.Cells.Item(140, 2).Value = 80
$ws.Cells.Item(140, 3).Value = 23
$ws.Cells.Item(140, 4).Value = 4
$ws.Cells.Item(140, 5).Value = 75

# Row 141 becomes Liechtenstein, carrying its previous (unchanged) figures
$ws.Cells.Item(141, 1).Value = "Liechtenstein"
$ws.Cells.Item(141, 2).Value = 79
$ws.Cells.Item(141, 4).Value = 55
$ws.Cells.Item(141, 5).Value = 23
$ws.Cells.Item(141, 8).Value = 1
